# Scheduled runner update: refresh market-board derived price/profit figures
# across the Ixion server profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Values below come from the latest Universalis price pull; only the
# price/profit columns (H:N) are touched, row identity columns (A:G) are
# left untouched.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# ALC
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("ALC")

# Row 21 - Book and a Hard Place / Engraved Hard Leather Grimoire
$ws.Range("H21").Value = 42707.6
$ws.Range("J21").Value = 26500
$ws.Range("L21").Value = 26500
$ws.Range("N21").Value = -27436

# Row 23 - There's Something about Bury / Hard Leather Grimoire
$ws.Range("H23").Value = 42707.6
$ws.Range("J23").Value = 26500
$ws.Range("L23").Value = 26500
$ws.Range("N23").Value = -26968

# Row 137 - Cutting Edge of Culinary Quality / Magnesia Whetstone
$ws.Range("H137").Value = 2117.6316
$ws.Range("I137").Value = 1323.3448
$ws.Range("J137").Value = 4677
$ws.Range("K137").Value = 3970.0344
$ws.Range("L137").Value = 14031
$ws.Range("M137").Value = -1420.0344
$ws.Range("N137").Value = -19131

# ---------------------------------------------------------------
# ARM
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("ARM")

# Row 13 - Get into Their Heads / Bronze Chain Coif
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("M13").ClearContents()

# Row 32 - Ingot We Trust / Steel Ingot
$ws.Range("H32").Value = 8532.299999999999
$ws.Range("I32").Value = 6187.633
$ws.Range("K32").Value = 6187.633
$ws.Range("M32").Value = -5900.633

# Row 45 - Hollow Hallmarks / Mythril Ingot
$ws.Range("H45").Value = 5245.231
$ws.Range("I45").Value = 5832.913
$ws.Range("J45").Value = 739.6667
$ws.Range("K45").Value = 5832.913
$ws.Range("L45").Value = 739.6667
$ws.Range("M45").Value = -5455.913
$ws.Range("N45").Value = -1493.6667

# Row 61 - Dealing with the Tough Stuff / Cobalt Ingot
$ws.Range("H61").Value = 196818.89
$ws.Range("I61").Value = 5540.1377
$ws.Range("K61").Value = 5540.1377
$ws.Range("M61").Value = -5328.1377

# Row 74 - As the Bolt Flies / Titanium Nugget
$ws.Range("H74").Value = 1642
$ws.Range("I74").Value = 1267.8667
$ws.Range("K74").Value = 1267.8667
$ws.Range("M74").Value = -393.8667

# Row 77 - Heavy Metal Banned (L) / Titanium Nugget
$ws.Range("H77").Value = 1642
$ws.Range("I77").Value = 1267.8667
$ws.Range("K77").Value = 6339.333500000001
$ws.Range("M77").Value = -1971.333500000001

# Row 132 - Don't Bore Me, Ore Me / Mountain Chromite Ingot
$ws.Range("H132").Value = 2781106
$ws.Range("I132").Value = 2269.0688
$ws.Range("K132").Value = 6807.2064
$ws.Range("M132").Value = -4277.2064

# Row 136 - Metal with Mettle / Cobalt Tungsten Ingot
$ws.Range("H136").Value = 196818.89
$ws.Range("I136").Value = 5540.1377
$ws.Range("K136").Value = 16620.4131
$ws.Range("M136").Value = -14070.4131

# ---------------------------------------------------------------
# BSM
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("BSM")

# Row 5 - Axe Me Anything / Bronze War Axe
$ws.Range("H5").Value = 3666.6667

# Row 16 - Port of Call: Ul'dah / Bronze Knuckles
$ws.Range("H16").Value = 4250
$ws.Range("I16").Value = 4250
$ws.Range("K16").Value = 4250
$ws.Range("M16").Value = -4080

# Row 21 - Awl or Nothing / Iron Awl
$ws.Range("H21").Value = 23500
$ws.Range("J21").Value = 23500
$ws.Range("L21").Value = 23500
$ws.Range("N21").Value = -23972

# Row 86 - Through Thick and Thin / Adamantite Nugget
$ws.Range("H86").Value = 1683.9565
$ws.Range("I86").Value = 1625.2858
$ws.Range("J86").Value = 2300
$ws.Range("K86").Value = 1625.2858
$ws.Range("L86").Value = 2300
$ws.Range("M86").Value = -502.2858000000001
$ws.Range("N86").Value = -4546

# Row 89 - Piercing Eyes Deserve Piercing Shafts (L) / Adamantite Nugget
$ws.Range("H89").Value = 1683.9565
$ws.Range("I89").Value = 1625.2858
$ws.Range("J89").Value = 2300
$ws.Range("K89").Value = 8126.429
$ws.Range("L89").Value = 11500
$ws.Range("M89").Value = -2510.429
$ws.Range("N89").Value = -22732

# Row 105 - Ingot to Wing It / Molybdenum Ingot
$ws.Range("H105").Value = 2591.6667
$ws.Range("I105").Value = 2700
$ws.Range("J105").Value = 2483.3333
$ws.Range("K105").Value = 2700
$ws.Range("L105").Value = 2483.3333
$ws.Range("M105").Value = -953
$ws.Range("N105").Value = -5977.3333

# ---------------------------------------------------------------
# CRP
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("CRP")

# Row 31 - Wall Not Found / Walnut Lumber
$ws.Range("H31").Value = 5058.069
$ws.Range("I31").Value = 1726.0769
$ws.Range("J31").Value = 7765.3125
$ws.Range("K31").Value = 1726.0769
$ws.Range("L31").Value = 7765.3125
$ws.Range("M31").Value = -1431.0769
$ws.Range("N31").Value = -8355.3125

# Row 34 - Armoires of the Rich and Famous / Walnut Lumber
$ws.Range("H34").Value = 5058.069
$ws.Range("I34").Value = 1726.0769
$ws.Range("J34").Value = 7765.3125
$ws.Range("K34").Value = 1726.0769
$ws.Range("L34").Value = 7765.3125
$ws.Range("M34").Value = -1524.0769
$ws.Range("N34").Value = -8169.3125

# Row 58 - You Do the Heavy Lifting / Mahogany Lumber
$ws.Range("H58").Value = 258169.84
$ws.Range("I58").Value = 1483.591
$ws.Range("J58").Value = 590352.0600000001
$ws.Range("K58").Value = 1483.591
$ws.Range("L58").Value = 590352.0600000001
$ws.Range("M58").Value = -1280.591
$ws.Range("N58").Value = -590758.0600000001

# Row 136 - Turali Quality / Dark Mahogany Lumber
$ws.Range("H136").Value = 258169.84
$ws.Range("I136").Value = 1483.591
$ws.Range("J136").Value = 590352.0600000001
$ws.Range("K136").Value = 4450.772999999999
$ws.Range("L136").Value = 1771056.18
$ws.Range("M136").Value = -1900.772999999999
$ws.Range("N136").Value = -1776156.18

# ---------------------------------------------------------------
# CUL
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("CUL")

# Row 5 - What a Sap / Maple Syrup
$ws.Range("H5").Value = 6860.8945
$ws.Range("I5").Value = 10425.8
$ws.Range("K5").Value = 31277.4
$ws.Range("M5").Value = -31165.4

# Row 131 - The Mountain Steeped / Tsai tou Vounou
$ws.Range("H131").Value = 2273819
$ws.Range("I131").Value = 7692847
$ws.Range("J131").Value = 1323.2903
$ws.Range("K131").Value = 23078541
$ws.Range("L131").Value = 3969.8709
$ws.Range("M131").Value = -23073501
$ws.Range("N131").Value = -14049.8709

# Row 135 - Not-so-secret Ingredient / Royal Maple Syrup
$ws.Range("H135").Value = 6860.8945
$ws.Range("I135").Value = 10425.8
$ws.Range("K135").Value = 93832.2
$ws.Range("M135").Value = -91297.2

# ---------------------------------------------------------------
# GSM
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("GSM")

# Row 80 - Needs More Prayerbell / Hardsilver Ingot
$ws.Range("H80").Value = 5435.375
$ws.Range("I80").Value = 7274.4
$ws.Range("J80").Value = 2370.3333
$ws.Range("K80").Value = 7274.4
$ws.Range("L80").Value = 2370.3333
$ws.Range("M80").Value = -6276.4
$ws.Range("N80").Value = -4366.3333

# Row 83 - With a Noise That Reaches Heaven (L) / Hardsilver Ingot
$ws.Range("H83").Value = 5435.375
$ws.Range("I83").Value = 7274.4
$ws.Range("J83").Value = 2370.3333
$ws.Range("K83").Value = 36372
$ws.Range("L83").Value = 11851.6665
$ws.Range("M83").Value = -31380
$ws.Range("N83").Value = -21835.6665

# ---------------------------------------------------------------
# LTW
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("LTW")

# Row 40 - Best Served Toad / Toad Leather
$ws.Range("H40").Value = 76926080
$ws.Range("I40").Value = 76926080
$ws.Range("K40").Value = 76926080
$ws.Range("M40").Value = -76925944

# Row 93 - Hide to Go Seek / Gagana Leather
$ws.Range("H93").Value = 1671.1428
$ws.Range("I93").Value = 1329.2
$ws.Range("J93").Value = 2526
$ws.Range("K93").Value = 1329.2
$ws.Range("L93").Value = 2526
$ws.Range("M93").Value = -81.20000000000005
$ws.Range("N93").Value = -5022

# Row 136 - Respect for Br'aax / Br'aax Leather
$ws.Range("H136").Value = 12074.296
$ws.Range("I136").Value = 9117.056
$ws.Range("J136").Value = 17988.777
$ws.Range("K136").Value = 27351.168
$ws.Range("L136").Value = 53966.33099999999
$ws.Range("M136").Value = -24801.168
$ws.Range("N136").Value = -59066.33099999999

# Row 140 - Worqor Zormor or Bust / Gargantuaskin Shoes of Healing
$ws.Range("H140").Value = 58583.855
$ws.Range("J140").Value = 58583.855
$ws.Range("L140").Value = 58583.855
$ws.Range("N140").Value = -68943.85500000001

# ---------------------------------------------------------------
# WVR
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("WVR")

# Row 15 - Workplace Safety / Cotton Scarf
$ws.Range("H15").Value = 7100
$ws.Range("J15").Value = 7100
$ws.Range("L15").Value = 7100
$ws.Range("N15").Value = -7676

# Row 132 - Comfy Cabins / Snow Cotton Cloth
$ws.Range("H132").Value = 1218.8983
$ws.Range("I132").Value = 950.13464
$ws.Range("K132").Value = 2850.40392
$ws.Range("M132").Value = -320.4039199999997

# Row 136 - Weaving the Envelope / Sarcenet Cloth
$ws.Range("H136").Value = 2449.9138
$ws.Range("I136").Value = 2634.3333
$ws.Range("J136").Value = 2252.3215
$ws.Range("K136").Value = 7902.999899999999
$ws.Range("L136").Value = 6756.9645
$ws.Range("M136").Value = -5352.999899999999
$ws.Range("N136").Value = -11856.9645
